$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows before the old "final" summary row (row 76) ---
$ws.Rows("76:79").Insert()

# --- Fill in the new rows (order chosen to match original authoring order
#     of the shared-string table: label first, then new_var names, then
#     the code formulas in 198,167,216,288 order) ---
$ws.Range("F76").Value = "Restriction with missing weight data =1, else = 0"

$ws.Range("B76").Value = "restrict_missingwt_noex.167"
$ws.Range("B77").Value = "restrict_missingwt_noex.198"
$ws.Range("B78").Value = "restrict_missingwt_noex.216"
$ws.Range("B79").Value = "restrict_missingwt_noex.288"

$ws.Range("K77").Value = "sub_restrict_noex.198 == 1 & is.na(bmiz_bestavail.192) ~ 1, sub_restrict_noex.198 == 1 & is.na(bmiz_drop.192) ~ 1, sub_restrict_noex.167 == 0 ~ 0, sub_restrict_noex.198 == 1 & !is.na(bmiz_bestavail.192) & !is.na(bmiz_drop.192) ~ 0"
$ws.Range("K76").Value = "sub_restrict_noex.167 == 1 & is.na(bmiz_bestavail.168) ~ 1, sub_restrict_noex.167 == 1 & is.na(bmiz_drop.168) ~ 1, sub_restrict_noex.167 == 0 ~ 0, sub_restrict_noex.167 == 1 & !is.na(bmiz_bestavail.168) & !is.na(bmiz_drop.168) ~ 0"
$ws.Range("K78").Value = "sub_restrict_noex.216 == 1 & is.na(bmiz_bestavail.216) ~ 1, sub_restrict_noex.216 == 1 & is.na(bmiz_drop.216) ~ 1, sub_restrict_noex.216 == 0 ~ 0, sub_restrict_noex.216 == 1 & !is.na(bmiz_bestavail.216) & !is.na(bmiz_drop.216) ~ 0"
$ws.Range("K79").Value = "sub_restrict_noex.288 == 1 & is.na(bmi_bestavail.288) ~ 1, sub_restrict_noex.288 == 1 & is.na(bmiz_drop.288) ~ 1, sub_restrict_noex.288 == 0 ~ 0, sub_restrict_noex.288 == 1 & !is.na(bmi_bestavail.288) & !is.na(bmiz_drop.288) ~ 0"

$ws.Range("D76:D79").Value = "case_when"
$ws.Range("E76:E79").Value = 5
$ws.Range("F77:F79").Value = "Restriction with missing weight data =1, else = 0"

# --- Update the formula in the (now shifted) final summary row 80 ---
$ws.Range("A80").Formula = '=TEXTJOIN(",", TRUE, A2,B32:B79)'

# --- Re-apply the existing sort over the (now larger) data range ---
$ws.Sort.SortFields.Clear()
$keyRange = $ws.Range("E33:E84")
$ws.Sort.SortFields.Add($keyRange, 0, 1) | Out-Null
$sortRange = $ws.Range("A2:K84")
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Update the view: scrolled position and active selection ---
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K79").Select()
